# Auto-generated edit script applying numeric corrections to the
# Chocobo_Profits workbook across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 889030.75
$ws.Range("I19").Value = 1333438.2
$ws.Range("J19").Value = 215.6
$ws.Range("K19").Value = 1333438.2
$ws.Range("L19").Value = 215.6
$ws.Range("M19").Value = -1333263.2
$ws.Range("N19").Value = -565.6

# Row 113
$ws.Range("H113").Value = 6918.737
$ws.Range("I113").Value = 3855.9
$ws.Range("J113").Value = 10321.889
$ws.Range("K113").Value = 3855.9
$ws.Range("L113").Value = 10321.889
$ws.Range("M113").Value = -601.9000000000001
$ws.Range("N113").Value = -16829.889

# Row 129
$ws.Range("H129").Value = 900.6222
$ws.Range("J129").Value = 963.325
$ws.Range("L129").Value = 2889.975
$ws.Range("N129").Value = -12889.975

# Row 132
$ws.Range("H132").Value = 77235720
$ws.Range("I132").Value = 91276300
$ws.Range("J132").Value = 12500
$ws.Range("K132").Value = 273828900
$ws.Range("L132").Value = 37500
$ws.Range("M132").Value = -273826370
$ws.Range("N132").Value = -42560

# Row 137
$ws.Range("H137").Value = 2872.3833
$ws.Range("I137").Value = 2203.8235
$ws.Range("J137").Value = 6660.8887
$ws.Range("K137").Value = 6611.470499999999
$ws.Range("L137").Value = 19982.6661
$ws.Range("M137").Value = -4061.470499999999
$ws.Range("N137").Value = -25082.6661

# Row 138
$ws.Range("H138").Value = 2370.8206
$ws.Range("I138").Value = 1328.2778
$ws.Range("J138").Value = 3264.4285
$ws.Range("K138").Value = 3984.8334
$ws.Range("L138").Value = 9793.2855
$ws.Range("M138").Value = 1155.1666
$ws.Range("N138").Value = -20073.2855

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 641.46875
$ws.Range("I2").Value = 619.04346
$ws.Range("J2").Value = 698.7778
$ws.Range("K2").Value = 619.04346
$ws.Range("L2").Value = 698.7778
$ws.Range("M2").Value = -506.04346
$ws.Range("N2").Value = -924.7778

# Row 74
$ws.Range("H74").Value = 3263.9
$ws.Range("I74").Value = 3639.2144
$ws.Range("K74").Value = 3639.2144
$ws.Range("M74").Value = -2765.2144

# Row 77
$ws.Range("H77").Value = 3263.9
$ws.Range("I77").Value = 3639.2144
$ws.Range("K77").Value = 18196.072
$ws.Range("M77").Value = -13828.072

# Row 80
$ws.Range("H80").Value = 38241.5
$ws.Range("J80").Value = 38241.5
$ws.Range("L80").Value = 38241.5
$ws.Range("N80").Value = -40237.5

# Row 83
$ws.Range("H83").Value = 38241.5
$ws.Range("J83").Value = 38241.5
$ws.Range("L83").Value = 114724.5
$ws.Range("N83").Value = -124708.5

# Row 116
$ws.Range("H116").Value = 641.46875
$ws.Range("I116").Value = 619.04346
$ws.Range("J116").Value = 698.7778
$ws.Range("K116").Value = 619.04346
$ws.Range("L116").Value = 698.7778
$ws.Range("M116").Value = 1674.95654
$ws.Range("N116").Value = -5286.7778

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 641.46875
$ws.Range("I3").Value = 619.04346
$ws.Range("J3").Value = 698.7778
$ws.Range("K3").Value = 619.04346
$ws.Range("L3").Value = 698.7778
$ws.Range("M3").Value = -505.04346
$ws.Range("N3").Value = -926.7778

# Row 134
$ws.Range("H134").Value = 2643.5208
$ws.Range("I134").Value = 1580.3572
$ws.Range("J134").Value = 10085.667
$ws.Range("K134").Value = 4741.071599999999
$ws.Range("L134").Value = 30257.001
$ws.Range("M134").Value = -2206.071599999999
$ws.Range("N134").Value = -35327.001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2948.4055
$ws.Range("I31").Value = 1029.5217
$ws.Range("J31").Value = 6100.857
$ws.Range("K31").Value = 1029.5217
$ws.Range("L31").Value = 6100.857
$ws.Range("M31").Value = -734.5217
$ws.Range("N31").Value = -6690.857

# Row 34
$ws.Range("H34").Value = 2948.4055
$ws.Range("I34").Value = 1029.5217
$ws.Range("J34").Value = 6100.857
$ws.Range("K34").Value = 1029.5217
$ws.Range("L34").Value = 6100.857
$ws.Range("M34").Value = -827.5217
$ws.Range("N34").Value = -6504.857

# Row 86
$ws.Range("H86").Value = 2462.077
$ws.Range("I86").Value = 1800
$ws.Range("J86").Value = 2756.3333
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 2756.3333
$ws.Range("M86").Value = -677
$ws.Range("N86").Value = -5002.3333

# Row 89
$ws.Range("H89").Value = 2462.077
$ws.Range("I89").Value = 1800
$ws.Range("J89").Value = 2756.3333
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 13781.6665
$ws.Range("M89").Value = -3384
$ws.Range("N89").Value = -25013.6665

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 512.96
$ws.Range("I14").Value = 512.96
$ws.Range("K14").Value = 1538.88
$ws.Range("M14").Value = -1365.88

# Row 40
$ws.Range("H40").Value = 437.625
$ws.Range("I40").Value = 143.14285
$ws.Range("J40").Value = 666.6667
$ws.Range("K40").Value = 572.5714
$ws.Range("L40").Value = 2666.6668
$ws.Range("M40").Value = -503.5714
$ws.Range("N40").Value = -2804.6668

# Row 129
$ws.Range("H129").Value = 2106.2856
$ws.Range("I129").Value = 2593.2666
$ws.Range("J129").Value = 888.8333
$ws.Range("K129").Value = 7779.7998
$ws.Range("L129").Value = 2666.4999
$ws.Range("M129").Value = -2779.7998
$ws.Range("N129").Value = -12666.4999

# Row 131
$ws.Range("H131").Value = 1011.431
$ws.Range("J131").Value = 935.8333
$ws.Range("L131").Value = 2807.4999
$ws.Range("N131").Value = -12887.4999

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 21066.309
$ws.Range("I43").Value = 2008.5
$ws.Range("J43").Value = 24531.363
$ws.Range("K43").Value = 2008.5
$ws.Range("L43").Value = 24531.363
$ws.Range("M43").Value = -1857.5
$ws.Range("N43").Value = -24833.363

# Row 46
$ws.Range("H46").Value = 32303.834
$ws.Range("J46").Value = 32303.834
$ws.Range("L46").Value = 32303.834
$ws.Range("N46").Value = -32615.834

# Row 57
$ws.Range("H57").Value = 36656.668
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 36656.668
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 36656.668
$ws.Range("N57").Value = -38296.668
$ws.Range("M57").ClearContents()

# Row 80
$ws.Range("H80").Value = 83336000
$ws.Range("I80").Value = 125002500
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 125002500
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -125001502
$ws.Range("N80").Value = -4996

# Row 83
$ws.Range("H83").Value = 83336000
$ws.Range("I83").Value = 125002500
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 625012500
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -625007508
$ws.Range("N83").Value = -24984

# Row 122
$ws.Range("H122").Value = 2303.6667
$ws.Range("I122").Value = 1866.2307
$ws.Range("J122").Value = 3928.4285
$ws.Range("K122").Value = 5598.6921
$ws.Range("L122").Value = 11785.2855
$ws.Range("M122").Value = -3148.6921
$ws.Range("N122").Value = -16685.2855

# Row 134
$ws.Range("H134").Value = 51890.9
$ws.Range("J134").Value = 51890.9
$ws.Range("L134").Value = 155672.7
$ws.Range("N134").Value = -160742.7

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3280.3
$ws.Range("I22").Value = 2225.25
$ws.Range("J22").Value = 3983.6667
$ws.Range("K22").Value = 2225.25
$ws.Range("L22").Value = 3983.6667
$ws.Range("M22").Value = -1930.25
$ws.Range("N22").Value = -4573.6667

# Row 27
$ws.Range("H27").Value = 3280.3
$ws.Range("I27").Value = 2225.25
$ws.Range("J27").Value = 3983.6667
$ws.Range("K27").Value = 2225.25
$ws.Range("L27").Value = 3983.6667
$ws.Range("M27").Value = -2118.25
$ws.Range("N27").Value = -4197.6667

# Row 135
$ws.Range("H135").Value = 50429
$ws.Range("J135").Value = 50429
$ws.Range("L135").Value = 50429
$ws.Range("N135").Value = -60569

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 329.5
$ws.Range("I113").Value = 324.9
$ws.Range("J113").Value = 335.25
$ws.Range("K113").Value = 974.6999999999999
$ws.Range("L113").Value = 1005.75
$ws.Range("M113").Value = 1195.3
$ws.Range("N113").Value = -5345.75

# Row 126
$ws.Range("H126").Value = 346092.06
$ws.Range("I126").Value = 2300.4211
$ws.Range("J126").Value = 890428.8
$ws.Range("K126").Value = 6901.263300000001
$ws.Range("L126").Value = 2671286.4
$ws.Range("M126").Value = -4431.263300000001
$ws.Range("N126").Value = -2676226.4
